# Update the model comparison metrics in the worksheet.
# The target values are strings that look numeric (currency amounts /
# large negative numbers), so force a Text number format on those cells
# first to ensure Excel keeps them as literal text instead of silently
# converting them into numeric cells (which would also lose precision
# on the long decimal values).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:C4").NumberFormat = "@"

$ws.Range("B2").Value = "€27,144.34"
$ws.Range("C2").Value = "€27,459.02"

$ws.Range("B3").Value = "-1650788376.1205"
$ws.Range("C3").Value = "-1689285588.1029"

$ws.Range("B4").Value = "€22,130.19"
$ws.Range("C4").Value = "€22,360.75"
